# "create account is completed"
#
# Adds a "plan" column to the GitHubSync sheet (sheet3), recording a FREE
# plan for the bhautik-vasebh / bh account that was just created.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GitHubSync")
$ws.Activate()

# New header + values. Write column C first (plan/FREE), then refresh the
# existing account columns (A/B) with the new account's username/company,
# matching the order the data was entered in the original session.
$ws.Range("C1").Value = "plan"
$ws.Range("C2").Value = "FREE"
$ws.Range("A2").Value = "bhautik-vasebh"
$ws.Range("B2").Value = "bh"

# Resize columns A and B to fit their (now-updated) contents.
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()

# Keep the page laid out in portrait orientation.
$ws.PageSetup.Orientation = 1

# Leave the selection where the user ended up after entering the data.
$ws.Range("C5").Select()
